$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-register the new item names in the shared-string table in the exact
# order they were first typed, using a scratch column far outside the used
# range. This fixes each name to a stable shared-string index before the
# G/I columns (which reference these names out of row order) are populated.
$newItemNames = @(
    "Icicle",
    "Bone Sword",
    "Obsidian Blade",
    "Claymore",
    "Weed Whacker",
    "Mist Saber",
    "Magma Blade",
    "Lightsaber",
    "Golden Khopesh",
    "Iron",
    "Copper",
    "Templar sword",
    "Destroyer",
    "Wooden Shaft",
    "Treasure Blade",
    "Ruby Blade",
    "Dragonbone Blade",
    "Elven Sword",
    "Serpent Sword",
    "Amethyst Shaft",
    "Captain's Sword",
    "Bronze Destroyer",
    "Iron Khopesh",
    "Opal",
    "Cleaver",
    "Magic Sword",
    "Ceremony Sword",
    "Relic Sword",
    "Eversteel Sword",
    "Jeweled sword",
    "Rhrodochrosite Blade",
    "Coral Sword",
    "Quartz Sword"
)
for ($i = 0; $i -lt $newItemNames.Count; $i++) {
    $ws.Cells.Item(1, 26 + $i).Value = $newItemNames[$i]
}
for ($i = 0; $i -lt $newItemNames.Count; $i++) {
    $ws.Cells.Item(1, 26 + $i).Value = ""
}

# Populate the new "item generation" lookup columns G (name) and I (weight)
# for rows 1-44.
$ws.Range("G1").Value = "Icicle"
$ws.Range("I1").Value = 20
$ws.Range("G2").Value = "Bone Sword"
$ws.Range("I2").Value = 2
$ws.Range("G3").Value = "Obsidian Blade"
$ws.Range("I3").Value = 16
$ws.Range("G4").Value = "Claymore"
$ws.Range("I4").Value = 15
$ws.Range("G5").Value = "Destroyer"
$ws.Range("I5").Value = 10
$ws.Range("G6").Value = "Weed Whacker"
$ws.Range("I6").Value = 4
$ws.Range("G7").Value = "Mist Saber"
$ws.Range("I7").Value = 26
$ws.Range("G8").Value = "Magma Blade"
$ws.Range("I8").Value = 25
$ws.Range("G9").Value = "Lightsaber"
$ws.Range("I9").Value = 23
$ws.Range("G10").Value = "Golden Khopesh"
$ws.Range("I10").Value = 11
$ws.Range("G11").Value = "Iron"
$ws.Range("I11").Value = 5
$ws.Range("G12").Value = "Steel"
$ws.Range("I12").Value = 9
$ws.Range("G13").Value = "Clot"
$ws.Range("I13").Value = 6
$ws.Range("G14").Value = "Ironleaf"
$ws.Range("I14").Value = 7
$ws.Range("G15").Value = "Copper"
$ws.Range("I15").Value = 3
$ws.Range("G16").Value = "Templar sword"
$ws.Range("I16").Value = 14
$ws.Range("G17").Value = "Quartz Sword"
$ws.Range("I17").Value = 12
$ws.Range("G18").Value = "Emberbronze"
$ws.Range("I18").Value = 24
$ws.Range("G19").Value = "Abyssal"
$ws.Range("I19").Value = 28
$ws.Range("G20").Value = "Opal"
$ws.Range("I20").Value = 15
$ws.Range("G21").Value = "Golden"
$ws.Range("I21").Value = 9
$ws.Range("G22").Value = "Lithium"
$ws.Range("I22").Value = 16
$ws.Range("G23").Value = "Wooden Shaft"
$ws.Range("I23").Value = 1
$ws.Range("G24").Value = "Treasure Blade"
$ws.Range("I24").Value = 26
$ws.Range("G25").Value = "Rhrodochrosite Blade"
$ws.Range("I25").Value = 23
$ws.Range("G26").Value = "Cleaver"
$ws.Range("I26").Value = 6
$ws.Range("G27").Value = "Serpent Sword"
$ws.Range("I27").Value = 27
$ws.Range("G28").Value = "Ceremony Sword"
$ws.Range("I28").Value = 30
$ws.Range("G29").Value = "Relic Sword"
$ws.Range("I29").Value = 29
$ws.Range("G30").Value = "Elven Sword"
$ws.Range("I30").Value = 10
$ws.Range("G31").Value = "Amethyst Shaft"
$ws.Range("I31").Value = 21
$ws.Range("G32").Value = "Captain's Sword"
$ws.Range("I32").Value = 24
$ws.Range("G33").Value = "Magic Sword"
$ws.Range("I33").Value = 18
$ws.Range("G34").Value = "Bronze Destroyer"
$ws.Range("I34").Value = 13
$ws.Range("G35").Value = "Iron Khopesh"
$ws.Range("I35").Value = 8
$ws.Range("G36").Value = "Dragonbone Blade"
$ws.Range("I36").Value = 20
$ws.Range("G37").Value = "Azure"
$ws.Range("I37").Value = 17
$ws.Range("G38").Value = "Coral Sword"
$ws.Range("I38").Value = 19
$ws.Range("G39").Value = "Magma Blade"
$ws.Range("I39").Value = 22
$ws.Range("G40").Value = "Eversteel Sword"
$ws.Range("I40").Value = 14
$ws.Range("G41").Value = "Jeweled sword"
$ws.Range("I41").Value = 24
$ws.Range("G42").Value = "Jeweled sword"
$ws.Range("I42").Value = 27
$ws.Range("G43").Value = "Jeweled sword"
$ws.Range("I43").Value = 25
$ws.Range("G44").Value = "Ruby Blade"
$ws.Range("I44").Value = 15

$ws.Range("I37").Select()

Write-Output "done"
